# Add a new participant e-mail address (kelker@kelker.de) below the
# existing one, mirroring how jtuttas@gmx.net is stored in A1: the cell
# holds the e-mail text and carries a "mailto:" hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "kelker@kelker.de"
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:kelker@kelker.de")

# Excel leaves the selection on the cell below the newly entered one
# (as if the user had pressed Enter after typing the address).
$ws.Range("A3").Select()
